$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain a text string
# (prevents Excel from auto-converting numeric-looking strings to numbers),
# then strips the temporary number-format override so the cell keeps its
# original (default) style, same as the source workbook.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "27.818.91"
$ws.Range("E2").Value = "  +0.71%  "
Set-TextValue "D3" "1.850.91"
$ws.Range("E3").Value = "  +0.11%  "
Set-TextValue "D5" "334.84"
Set-TextValue "D6" "1.005"
$ws.Range("E6").Value = "  -0.25%  "
Set-TextValue "D7" "0.4650"
$ws.Range("E7").Value = "  +0.96%  "
Set-TextValue "D8" "0.3868"
$ws.Range("E8").Value = "  -0.53%  "
Set-TextValue "D9" "46.67"
$ws.Range("E9").Value = "  +1.46%  "
Set-TextValue "D10" "0.07915"
$ws.Range("E10").Value = "  -0.38%  "
Set-TextValue "D11" "0.9698"
$ws.Range("E11").Value = "  -3.46%  "
Set-TextValue "D12" "21.34"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.843.00"
$ws.Range("E13").Value = "  -1.03%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D14" "5.895"
$ws.Range("E14").Value = "  -1.32%  "
Set-TextValue "D15" "7.170"
$ws.Range("E15").Value = "  -0.08%  "
$ws.Range("E16").Value = "  -0.37%  "
Set-TextValue "D17" "90.12"
$ws.Range("E17").Value = "  +1.89%  "
Set-TextValue "D18" "0.06610"
$ws.Range("E18").Value = "  -1.13%  "
Set-TextValue "D19" "0.00001029"
$ws.Range("E19").Value = "  -0.80%  "
Set-TextValue "D20" "17.38"
$ws.Range("E20").Value = "  +0.81%  "
Set-TextValue "D21" "1.006"
$ws.Range("E21").Value = "  -0.10%  "
Set-TextValue "D22" "27.806.25"
$ws.Range("E22").Value = "  +0.58%  "
Set-TextValue "D23" "5.349"
$ws.Range("E23").Value = "  -1.18%  "
Set-TextValue "D24" "10.86"
$ws.Range("E24").Value = "  -1.17%  "
Set-TextValue "D25" "2.294"
$ws.Range("E25").Value = "  -0.84%  "
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue "D26" "2.084.16"
$ws.Range("E26").Value = "  +0.21%  "
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D27" "158.64"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  -0.31%  "
Set-TextValue "D29" "2.071"
$ws.Range("E29").Value = "  -2.61%  "
Set-TextValue "D30" "5.368"
$ws.Range("E30").Value = "  -1.50%  "
Set-TextValue "D31" "118.76"
$ws.Range("E31").Value = "  -1.72%  "
Set-TextValue "D32" "0.09416"
$ws.Range("E32").Value = "  +0.05%  "
Set-TextValue "D33" "0.9480"
$ws.Range("E33").Value = "  -3.07%  "
Set-TextValue "D34" "3.588"
$ws.Range("E34").Value = "  -0.74%  "
Set-TextValue "D35" "5.265"
$ws.Range("E35").Value = "  -0.79%  "
Set-TextValue "D36" "1.329"
$ws.Range("E36").Value = "  -0.96%  "
Set-TextValue "D37" "0.06017"
$ws.Range("E37").Value = "  -0.13%  "
Set-TextValue "D38" "0.02218"
$ws.Range("E38").Value = "  -0.68%  "
Set-TextValue "D39" "8.247"
$ws.Range("E39").Value = "  -0.90%  "
$ws.Range("E40").Value = "  -0.26%  "
Set-TextValue "D41" "1.156"
$ws.Range("E41").Value = "  -2.43%  "
Set-TextValue "D42" "0.5825"
$ws.Range("E42").Value = "  -1.67%  "
Set-TextValue "D43" "0.1845"
$ws.Range("E43").Value = "  -1.32%  "
Set-TextValue "D44" "10.08"
$ws.Range("E44").Value = "  -2.71%  "
Set-TextValue "D45" "1.280"
$ws.Range("E45").Value = "  +2.93%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D46" "0.5457"
$ws.Range("E46").Value = "  -2.37%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "11.92"
$ws.Range("E47").Value = "  -2.02%  "
Set-TextValue "D48" "1.932"
$ws.Range("E48").Value = "  +0.94%  "
Set-TextValue "D49" "0.06846"
$ws.Range("E49").Value = "  +2.12%  "
$ws.Range("B50").Value = "PaxosStandard"
$ws.Range("C50").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
Set-TextValue "D50" "1.049"
$ws.Range("E50").Value = "  -29.70%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D51" "110.93"
$ws.Range("E51").Value = "  -0.03%  "
